$wb = $excel.ActiveWorkbook

# Sheet 1: "Заказы" (Orders) - update latitude (B) / longitude (C) values, rows 2-46
$wsOrders = $wb.Worksheets.Item(1)
$wsOrders.Range("B2").Value = 55.68444453800417
$wsOrders.Range("C2").Value = 37.57553578589383
$wsOrders.Range("B3").Value = 55.67948788095391
$wsOrders.Range("C3").Value = 37.73522421131663
$wsOrders.Range("B4").Value = 55.70251682158581
$wsOrders.Range("C4").Value = 37.59679890499072
$wsOrders.Range("B5").Value = 55.67450824280149
$wsOrders.Range("C5").Value = 37.67161683292756
$wsOrders.Range("B6").Value = 55.71201203249306
$wsOrders.Range("C6").Value = 37.56095890833359
$wsOrders.Range("B7").Value = 55.71002716873189
$wsOrders.Range("C7").Value = 37.57445901435683
$wsOrders.Range("B8").Value = 55.74307564965378
$wsOrders.Range("C8").Value = 37.67353730835671
$wsOrders.Range("B9").Value = 55.70794586738244
$wsOrders.Range("C9").Value = 37.5155211703763
$wsOrders.Range("B10").Value = 55.70566427688573
$wsOrders.Range("C10").Value = 37.63377173727263
$wsOrders.Range("B11").Value = 55.71446251453077
$wsOrders.Range("C11").Value = 37.73127929139365
$wsOrders.Range("B12").Value = 55.66529473739025
$wsOrders.Range("C12").Value = 37.65875383503119
$wsOrders.Range("B13").Value = 55.76650305010796
$wsOrders.Range("C13").Value = 37.67643968207631
$wsOrders.Range("B14").Value = 55.69511294207584
$wsOrders.Range("C14").Value = 37.73507280763371
$wsOrders.Range("B15").Value = 55.66120071251704
$wsOrders.Range("C15").Value = 37.59824071791409
$wsOrders.Range("B16").Value = 55.71691655903573
$wsOrders.Range("C16").Value = 37.49599886319182
$wsOrders.Range("B17").Value = 55.78321235662322
$wsOrders.Range("C17").Value = 37.55102611159206
$wsOrders.Range("B18").Value = 55.75955603165747
$wsOrders.Range("C18").Value = 37.53418321376897
$wsOrders.Range("B19").Value = 55.79916183382227
$wsOrders.Range("C19").Value = 37.61070137759584
$wsOrders.Range("B20").Value = 55.74643952086186
$wsOrders.Range("C20").Value = 37.50341270110896
$wsOrders.Range("B21").Value = 55.8193043852534
$wsOrders.Range("C21").Value = 37.54766770218015
$wsOrders.Range("B22").Value = 55.75559622911472
$wsOrders.Range("C22").Value = 37.63597152864205
$wsOrders.Range("B23").Value = 55.67665421830194
$wsOrders.Range("C23").Value = 37.68190332216877
$wsOrders.Range("B24").Value = 55.72252272070857
$wsOrders.Range("C24").Value = 37.7341602222447
$wsOrders.Range("B25").Value = 55.73990078453081
$wsOrders.Range("C25").Value = 37.7293620568059
$wsOrders.Range("B26").Value = 55.79835675044812
$wsOrders.Range("C26").Value = 37.45008770291449
$wsOrders.Range("B27").Value = 55.65125928042025
$wsOrders.Range("C27").Value = 37.61993628386446
$wsOrders.Range("B28").Value = 55.71478725332855
$wsOrders.Range("C28").Value = 37.5751873780851
$wsOrders.Range("B29").Value = 55.78064293979297
$wsOrders.Range("C29").Value = 37.55606655974139
$wsOrders.Range("B30").Value = 55.68606143474795
$wsOrders.Range("C30").Value = 37.52943375936608
$wsOrders.Range("B31").Value = 55.69626062163972
$wsOrders.Range("C31").Value = 37.6876735666027
$wsOrders.Range("B32").Value = 55.75604903038228
$wsOrders.Range("C32").Value = 37.47946098310879
$wsOrders.Range("B33").Value = 55.75183617624158
$wsOrders.Range("C33").Value = 37.56891589443138
$wsOrders.Range("B34").Value = 55.70240577242451
$wsOrders.Range("C34").Value = 37.6712653294931
$wsOrders.Range("B35").Value = 55.75547556771942
$wsOrders.Range("C35").Value = 37.7326695374279
$wsOrders.Range("B36").Value = 55.71823477831779
$wsOrders.Range("C36").Value = 37.46397928098686
$wsOrders.Range("B37").Value = 55.77819091680744
$wsOrders.Range("C37").Value = 37.63454577765092
$wsOrders.Range("B38").Value = 55.66507806523065
$wsOrders.Range("C38").Value = 37.70267233299501
$wsOrders.Range("B39").Value = 55.68996882645116
$wsOrders.Range("C39").Value = 37.46856778456943
$wsOrders.Range("B40").Value = 55.79546550716621
$wsOrders.Range("C40").Value = 37.55812694306487
$wsOrders.Range("B41").Value = 55.74469583208078
$wsOrders.Range("C41").Value = 37.74482413564305
$wsOrders.Range("B42").Value = 55.77060872411177
$wsOrders.Range("C42").Value = 37.63068928072016
$wsOrders.Range("B43").Value = 55.76808816408723
$wsOrders.Range("C43").Value = 37.53883926558113
$wsOrders.Range("B44").Value = 55.67956458238026
$wsOrders.Range("C44").Value = 37.58811896900087
$wsOrders.Range("B45").Value = 55.71076548708987
$wsOrders.Range("C45").Value = 37.68635213515589
$wsOrders.Range("B46").Value = 55.74793554571764
$wsOrders.Range("C46").Value = 37.69384468198596

# Sheet 2: "Курьеры" (Couriers) - update work schedule (C) / profile (D) text values
$wsCouriers = $wb.Worksheets.Item(2)
$wsCouriers.Range("C2").Value = "(2020-10-01 08:00:00 - 2020-10-01 22:00:00)"
$wsCouriers.Range("D3").Value = "car"
$wsCouriers.Range("C4").Value = "(2020-10-01 08:00:00 - 2020-10-01 19:00:00)"
$wsCouriers.Range("D5").Value = "car"
$wsCouriers.Range("C6").Value = "(2020-10-01 12:00:00 - 2020-10-01 19:00:00)"

# Sheet 3: "Склады" (Depots) - update latitude (C) / longitude (D) / work schedule (E)
$wsDepots = $wb.Worksheets.Item(3)
$wsDepots.Range("C2").Value = 55.75530405177981
$wsDepots.Range("D2").Value = 37.46851600393571
$wsDepots.Range("E2").Value = "(2020-10-01 09:00:00 - 2020-10-01 22:00:00)"
$wsDepots.Range("C3").Value = 55.68947874064456
$wsDepots.Range("D3").Value = 37.52911383980694
$wsDepots.Range("E3").Value = "(2020-10-01 11:00:00 - 2020-10-01 19:00:00)"
$wsDepots.Range("C4").Value = 55.7771001718272
$wsDepots.Range("D4").Value = 37.54383200505681
$wsDepots.Range("E4").Value = "(2020-10-01 10:00:00 - 2020-10-01 21:00:00)"

# Sheet 4: "Профили" (Profiles) - rename "bicycle"/"foot" profiles to "car"
$wsProfiles = $wb.Worksheets.Item(4)
$wsProfiles.Range("B3").Value = "car"
$wsProfiles.Range("C3").Value = "car"
$wsProfiles.Range("B5").Value = "car"
$wsProfiles.Range("C5").Value = "car"
